$p = $ppt.ActivePresentation

# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / "Office" colours (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colours (used by the Slide Master)
# The target state swaps which colour set backs the slide design: the slide
# master's effective theme becomes the plain "Office" colour scheme while the
# "Red Violet" colours move off of it. We reach the live (slide-master-backed)
# theme through ThemeColorScheme and overwrite every colour slot in place.

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# (RGB() isn't available in this host, so colours are passed as the packed
# 0x00BBGGRR long values PowerPoint's own RGB property getter/setter uses.)
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
